$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("225:226").Insert()

$ws.Range("A225").Value = 8
$ws.Range("B225").Value = "Terminal La Palmera de La Serena"
$ws.Range("C225").Value = "Coquimbo"
$ws.Range("D225").Value = 44694
$ws.Range("E225").Value = 4
$ws.Range("F225").Value = 100112003
$ws.Range("G225").Value = "Ajo"
$ws.Range("H225").Value = "Chino"
$ws.Range("I225").Value = "Primera"
$ws.Range("J225").Value = 540
$ws.Range("K225").Value = 19000
$ws.Range("L225").Value = 20000
$ws.Range("M225").Value = 19500
$ws.Range("N225").Value = "$/caja 10 kilos"
$ws.Range("O225").Value = "China"
$ws.Range("P225").Value = 1950
$ws.Range("Q225").Value = 10
$ws.Range("R225").Value = "Hortaliza"

$ws.Range("A226").Value = 8
$ws.Range("B226").Value = "Terminal La Palmera de La Serena"
$ws.Range("C226").Value = "Coquimbo"
$ws.Range("D226").Value = 44694
$ws.Range("E226").Value = 4
$ws.Range("F226").Value = 100112003
$ws.Range("G226").Value = "Ajo"
$ws.Range("H226").Value = "Chino"
$ws.Range("I226").Value = "Primera"
$ws.Range("J226").Value = 540
$ws.Range("K226").Value = 20000
$ws.Range("L226").Value = 21000
$ws.Range("M226").Value = 20500
$ws.Range("N226").Value = "$/malla 10 kilos"
$ws.Range("O226").Value = "China"
$ws.Range("P226").Value = 2050
$ws.Range("Q226").Value = 10
$ws.Range("R226").Value = "Hortaliza"
